# Loan RBI, Variable Instalments
#
# Insert a new (blank) column N on the "Repayment schedule" sheet,
# pushing the existing "Late" / paid-date-heading / "Outstanding"
# columns one slot to the right (N->O, O->P, P->Q), and make the
# "Repayment schedule" sheet the active / selected sheet+cell.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column before the old column N (Late).
$null = $wsRepay.Columns("N:N").Insert()

# New column gets a plain (non bestFit) width of 11 characters.
$wsRepay.Columns("N:N").ColumnWidth = 10.1

# Make "Repayment schedule" the active sheet/tab and move the selection
# to S8 (this also updates workbookView.activeTab and tabSelected on the
# sheetViews of both this sheet and the previously-active one).
$null = $wsRepay.Activate()
$null = $wsRepay.Range("S8").Select()

Write-Host "done"
